# Natmi following Dr Hou advice:
# recompute the existing Wnt2-Fzd8 rows (ligand/receptor-expressing cell
# counts + derived specificities) and add two new target-cluster rows
# (M1, M2) for the same ligand/receptor pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Wnt2"
$ws.Cells.Item(2,3).Value = "Fzd8"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.3884013333333334
$ws.Cells.Item(2,8).Value = 1.165204
$ws.Cells.Item(2,9).Value = 1
$ws.Cells.Item(2,10).Value = 1
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.062569
$ws.Cells.Item(2,14).Value = 9.187707
$ws.Cells.Item(2,15).Value = 0.2460077391286943
$ws.Cells.Item(2,16).Value = 0.2690593624267
$ws.Cells.Item(2,17).Value = 1.189505883025333
$ws.Cells.Item(2,18).Value = 10.705552947228
$ws.Cells.Item(2,19).Value = 0.2460077391286943
$ws.Cells.Item(2,20).Value = 0.2690593624267

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Wnt2"
$ws.Cells.Item(3,3).Value = "Fzd8"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.3884013333333334
$ws.Cells.Item(3,8).Value = 1.165204
$ws.Cells.Item(3,9).Value = 1
$ws.Cells.Item(3,10).Value = 1
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 6.147102
$ws.Cells.Item(3,14).Value = 18.441306
$ws.Cells.Item(3,15).Value = 0.4937797859292232
$ws.Cells.Item(3,16).Value = 0.5400483531609875
$ws.Cells.Item(3,17).Value = 2.387542612936
$ws.Cells.Item(3,18).Value = 21.487883516424
$ws.Cells.Item(3,19).Value = 0.4937797859292232
$ws.Cells.Item(3,20).Value = 0.5400483531609875

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Wnt2"
$ws.Cells.Item(4,3).Value = "Fzd8"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.3884013333333334
$ws.Cells.Item(4,8).Value = 1.165204
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,10).Value = 1
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.01759833333333333
$ws.Cells.Item(4,14).Value = 0.052795
$ws.Cells.Item(4,15).Value = 0.001413625683459368
$ws.Cells.Item(4,16).Value = 0.001546086421706485
$ws.Cells.Item(4,17).Value = 0.006835216131111112
$ws.Cells.Item(4,18).Value = 0.06151694518000001
$ws.Cells.Item(4,19).Value = 0.001413625683459368
$ws.Cells.Item(4,20).Value = 0.001546086421706485

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Wnt2"
$ws.Cells.Item(5,3).Value = "Fzd8"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.3884013333333334
$ws.Cells.Item(5,8).Value = 1.165204
$ws.Cells.Item(5,9).Value = 1
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.022088
$ws.Cells.Item(5,14).Value = 0.066264
$ws.Cells.Item(5,15).Value = 0.001774268250568265
$ws.Cells.Item(5,16).Value = 0.00194052222081558
$ws.Cells.Item(5,17).Value = 0.008579008650666667
$ws.Cells.Item(5,18).Value = 0.07721107785600001
$ws.Cells.Item(5,19).Value = 0.001774268250568265
$ws.Cells.Item(5,20).Value = 0.00194052222081558

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Wnt2"
$ws.Cells.Item(6,3).Value = "Fzd8"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.3884013333333334
$ws.Cells.Item(6,8).Value = 1.165204
$ws.Cells.Item(6,9).Value = 1
$ws.Cells.Item(6,10).Value = 1
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.1997185
$ws.Cells.Item(6,14).Value = 6.399437000000001
$ws.Cells.Item(6,15).Value = 0.2570245810080548
$ws.Cells.Item(6,16).Value = 0.1874056757697904
$ws.Cells.Item(6,17).Value = 1.242774931691334
$ws.Cells.Item(6,18).Value = 7.456649590148001
$ws.Cells.Item(6,19).Value = 0.2570245810080548
$ws.Cells.Item(6,20).Value = 0.1874056757697904
